$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly driver report update for 2025-04-19
# The "Good Drivers" table (rows 12-17) is refreshed: Total Samples counts
# are updated and the rows are re-sorted so that 21.60.2.1 and 22.50.1.1
# move to the top (their Driver Vintage date is not yet known / blank),
# while the remaining four drivers shift down, keeping their own
# Total Samples / % / Vintage values attached to their name.

# Force a literal text type (matches source data, where Driver Vintage is
# stored as plain text like "2024-11-10", not a real Excel date) for the
# rows whose vintage date text is actually changing, before writing the
# values - this stops Excel from auto-converting the date-look-alike
# strings into date serial numbers. E17's vintage text is unchanged by
# this update, so it is left alone.
$ws.Range("E14:E16").NumberFormat = "@"

$ws.Cells.Item(12, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Cells.Item(12, 2).Value = 56018
$ws.Cells.Item(12, 4).Value = 100
$ws.Cells.Item(12, 5).Value = $null

$ws.Cells.Item(13, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Cells.Item(13, 2).Value = 34244
$ws.Cells.Item(13, 4).Value = 100
$ws.Cells.Item(13, 5).Value = $null

$ws.Cells.Item(14, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Cells.Item(14, 2).Value = 442178
$ws.Cells.Item(14, 4).Value = 99.90000000000001
$ws.Cells.Item(14, 5).Value = "2024-11-10"

$ws.Cells.Item(15, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Cells.Item(15, 2).Value = 77849
$ws.Cells.Item(15, 4).Value = 99.90000000000001
$ws.Cells.Item(15, 5).Value = "2021-08-18"

$ws.Cells.Item(16, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Cells.Item(16, 2).Value = 59673
$ws.Cells.Item(16, 4).Value = 100
$ws.Cells.Item(16, 5).Value = "2020-08-05"

$ws.Cells.Item(17, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Cells.Item(17, 2).Value = 113652
# D17 (100) and E17 ("2019-12-14") are unchanged by this update - leave as-is.
